# Applies the edits described by the commit:
#  - Two name corrections in the schedule table (B6, B7)
#  - Removal of leftover placeholder/test text from several "helper" cells
#    (B14, B15, B16, B35, B36, B38) so they go back to being blank, matching
#    their sibling rows elsewhere in the workbook
#  - Recolour of the week-header band fill (was a near-black 282a36,
#    becomes a brown 93684c)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Name corrections -------------------------------------------------
$ws.Range("B6").Value = "Ben Corran"
$ws.Range("B7").Value = "Graham Haynes"

# --- 2. Clear stray placeholder text --------------------------------------
# B15 and B35 ("asdf" / "dddd") previously carried a bespoke bold/centred
# style (fontId 6). Their sibling "(x min)" rows elsewhere (e.g. B37, B56,
# B104 ...) are plain, unstyled cells, so after clearing the text we copy
# that plain look across via a formats-only paste (this reuses the already
# existing plain style instead of inventing a new one).
$ws.Range("B15").ClearContents() | Out-Null
$ws.Range("B35").ClearContents() | Out-Null

$ws.Range("B37").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$ws.Range("B56").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null

# B14, B16, B36, B38 ("asdf" / "basdf" / "asdf" / "David Moore") previously
# carried the same bespoke bold/centred style too. Their sibling spacer rows
# (e.g. B57, B59 ...) already use the plain grey-fill style with no text, so
# clear the content and copy that look across the same way.
$ws.Range("B14").ClearContents() | Out-Null
$ws.Range("B16").ClearContents() | Out-Null
$ws.Range("B36").ClearContents() | Out-Null
$ws.Range("B38").ClearContents() | Out-Null

$ws.Range("B57").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

$ws.Range("B59").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null

$ws.Range("B61").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null

$ws.Range("B80").Copy() | Out-Null
$ws.Range("B38").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 3. Recolour the week-header band fill --------------------------------
# RGB(147,104,76) == 0x93684C, packed as BGR for the OLE Color value.
$newBandColor = 147 + (104 * 256) + (76 * 65536)

$bandCellsRow1 = @("A2", "A24", "A46", "A69", "A94", "A117", "A141", "A165")
$bandCellsRow2 = @("B2", "B24", "B46", "B69", "B94", "B117", "B141", "B165")
$bandCellsRow3 = @("C2", "C24", "C46", "C69", "C94", "C117", "C141", "C165")
$bandCellsRow4 = @("D2", "D24", "D46", "D69", "D94", "D117", "D141", "D165")

foreach ($addr in $bandCellsRow1) { $ws.Range($addr).Interior.Color = $newBandColor }
foreach ($addr in $bandCellsRow2) { $ws.Range($addr).Interior.Color = $newBandColor }
foreach ($addr in $bandCellsRow3) { $ws.Range($addr).Interior.Color = $newBandColor }
foreach ($addr in $bandCellsRow4) { $ws.Range($addr).Interior.Color = $newBandColor }
